$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-7 (2000年, 2005年, 2006年, 2007年, 2008年, 2009年)
# This shifts remaining rows (2010-2013, currently rows 8-11) up to rows 2-5
$ws.Range("A2:C7").EntireRow.Delete()
